# Update "想去人数" (F) and "最低票价" (G) figures across sheets to match
# the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 12298
$ws1.Range("F3").Value  = 6915
$ws1.Range("F6").Value  = 438
$ws1.Range("F10").Value = 960
$ws1.Range("G10").Value = 108
$ws1.Range("F13").Value = 968
$ws1.Range("F14").Value = 3698
$ws1.Range("F16").Value = 998
$ws1.Range("F19").Value = 339
$ws1.Range("F21").Value = 250
$ws1.Range("F24").Value = 335
$ws1.Range("F25").Value = 5119
$ws1.Range("F27").Value = 1332
$ws1.Range("F28").Value = 266
$ws1.Range("F29").Value = 793
$ws1.Range("F30").Value = 1278

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 3722

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 535

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 535
$ws4.Range("F5").Value  = 12298
$ws4.Range("F6").Value  = 6915
$ws4.Range("F8").Value  = 3722
$ws4.Range("F11").Value = 438
$ws4.Range("F15").Value = 960
$ws4.Range("G15").Value = 108
$ws4.Range("F18").Value = 968
$ws4.Range("F19").Value = 3698
$ws4.Range("F21").Value = 998
$ws4.Range("F24").Value = 339
$ws4.Range("F26").Value = 250
$ws4.Range("F32").Value = 335
$ws4.Range("F33").Value = 5119
$ws4.Range("F35").Value = 1332
$ws4.Range("F38").Value = 266
$ws4.Range("F40").Value = 793
$ws4.Range("F41").Value = 1278
